# Daily refresh of the cryptos worksheet (price + 1h volume figures),
# plus the reshuffle at the bottom of the table where a new coin
# (BabyDogeCoin) enters the ranked list at row 46, pushing the rows below
# it down by one and dropping the previous last entry (Aptos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-45: coin identity (columns B/C) is unchanged; only the latest
# price (D) and/or 1h volume change (E) are refreshed.
$ws.Range("D2").Value = "'29.378.06"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "'1.881.82"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'0.7123"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'242.34"
$ws.Range("D8").Value = "'0.08048"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("D9").Value = "'0.3124"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'25.28"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "'0.08327"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "'1.899.70"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "'5.245"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'0.7193"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "'93.63"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "'6.323"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("D17").Value = "'0.000008541"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").Value = "'29.387.83"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'241.82"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'2.136.37"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'13.24"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.848"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").Value = "'164.00"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'9.062"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "'1.508"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'4.416"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "'4.343"
$ws.Range("D32").Value = "'1.198"
$ws.Range("E32").Value = "  -6.16%  "
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'1.183"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'0.7490"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'2.697"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "'0.01887"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "'1.288.08"
$ws.Range("E39").Value = "  +9.44%  "
$ws.Range("D40").Value = "'2.747"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "'6.603"
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("D42").Value = "'0.9196"
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("D43").Value = "'112.14"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("D44").Value = "'74.43"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  +0.07%  "

# Rows 46-51: BabyDogeCoin is newly listed at row 46; RocketPoolETH,
# RenderToken, Mantle, EnergySwap and TheSandbox each shift down one row,
# and Aptos (previously row 51) drops off the bottom of the table.
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000128"
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "'2.038.68"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.806"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5222"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.517"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4385"
$ws.Range("E51").Value = "  +1.82%  "
